# This workbook holds weekly price-report rows for
# "Agrícola del Norte S.A. de Arica" / "Zapallo italiano" / "Huracán".
# A new week of data (serial date 45077) is being added at the top of the
# data block that starts at row 486. Inserting two rows there shifts every
# existing row down by two (dimension grows from R551 to R553) while
# preserving each row's original values - exactly what the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 486, pushing rows 486:551 down
# to 488:553.
$ws.Rows.Item(486).Insert()
$ws.Rows.Item(486).Insert()

# --- New row 486: "Primera" quality for the new week ---
$ws.Range("A486").Value = 1
$ws.Range("B486").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C486").Value = "Arica y Parinacota"
$ws.Range("D486").Value = 45077
$ws.Range("E486").Value = 15
$ws.Range("F486").Value = 100112032
$ws.Range("G486").Value = "Zapallo italiano"
$ws.Range("H486").Value = "Huracán"
$ws.Range("I486").Value = "Primera"
$ws.Range("J486").Value = 120
$ws.Range("K486").Value = 6000
$ws.Range("L486").Value = 7000
$ws.Range("M486").Value = 6500
$ws.Range("N486").Value = "$/caja 70 unidades"
$ws.Range("O486").Value = "Región de Arica y Parinacota"
$ws.Range("P486").Value = 93
$ws.Range("Q486").Value = 70
$ws.Range("R486").Value = "Hortaliza"

# --- New row 487: "Segunda" quality for the new week ---
$ws.Range("A487").Value = 1
$ws.Range("B487").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C487").Value = "Arica y Parinacota"
$ws.Range("D487").Value = 45077
$ws.Range("E487").Value = 15
$ws.Range("F487").Value = 100112032
$ws.Range("G487").Value = "Zapallo italiano"
$ws.Range("H487").Value = "Huracán"
$ws.Range("I487").Value = "Segunda"
$ws.Range("J487").Value = 140
$ws.Range("K487").Value = 4000
$ws.Range("L487").Value = 5000
$ws.Range("M487").Value = 4500
$ws.Range("N487").Value = "$/caja 100 unidades"
$ws.Range("O487").Value = "Región de Arica y Parinacota"
$ws.Range("P487").Value = 45
$ws.Range("Q487").Value = 100
$ws.Range("R487").Value = "Hortaliza"
